# Update the "Data" worksheet of the VRU Headform (C) color-wise data
# workbook: widen column A, rename/re-sort the existing car entries and
# append the new car rows with their color-distribution data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A gets a little wider (20 -> 23 characters).
$ws.Columns.Item(1).ColumnWidth = 22.15

# Full row data (Car Name, D Green, Green, Yellow, Orange, Brown, Red,
# Default Red, Blue, Predicted headform score) for rows 2-15.
$rows = @{
    2  = @("BYD SEAL 2023 ",        0,                  3.72,  42.56, 26.86, 4.13,  16.94, 5.79, 0, 100)
    3  = @("XPENG G9 2023 ",        0,                  11.51, 29.37, 22.22, 10.32, 24.21, 2.38, 0, 100)
    4  = @("VinFast VF8 2023 ",     0,                  40.08, 20.25, 9.09,  4.13,  21.49, 4.96, 0, 100)
    5  = @("Honda ZR V 2023 ",      9.130000000000001,  27.83, 34.35, 8.26,  6.96,  8.26,  5.22, 0, 100)
    6  = @("BYD SEAL U 2023 ",      0,                  22.18, 33.47, 10.08, 13.71, 17.34, 3.23, 0, 100)
    7  = @("Volkswagen ID.7 2023 ",4.2,                 34.03, 43.28, 3.36,  2.52,  6.72,  5.88, 0, 100)
    8  = @("BMW 5 Series 2023 ",    0,                  62.82, 14.53, 5.98,  4.27,  8.970000000000001, 3.42, 0, 100)
    9  = @("smart #3 ",             7.33,               36.64, 23.71, 11.64, 5.17,  9.48,  6.03, 0, 100)
    10 = @("BYD Tang 2023 ",        0,                  21.37, 34.19, 12.39, 9.83,  17.09, 5.13, 0, 100)
    11 = @("Hyundai KONA 2023 ",    0,                  34.35, 35.22, 11.74, 3.91,  9.57,  5.22, 0, 100)
    12 = @("Kia EV9 2023 ",         0,                  38.31, 30.24, 10.48, 4.84,  12.1,  4.03, 0, 100)
    13 = @("NIO ET5 2023 ",         0,                  19.77, 39.53, 14.34, 6.59,  13.57, 6.2,  0, 100)
    14 = @("NIO EL7 2023 ",         0,                  27.91, 25.19, 17.83, 6.98,  18.99, 3.1,  0, 100)
    15 = @("Lexus RZ 2023 ",        0,                  40,    25.65, 14.78, 9.57,  4.78,  5.22, 0, 100)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
